# Trade #13 closed at 2026-02-17 07:53:50 - unknown UNKNOWN +0.000%
# Applies the new closed trade row to the "All Trades" and "MarketMaking"
# sheets, and updates the rolled-up stats on "Summary" and "Strategy Status".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.96
$summary.Range("B4").Value = -0.04
$summary.Range("B5").Value = -0.06
$summary.Range("B6").Value = 13
$summary.Range("B8").Value = 7
$summary.Range("B9").Value = 38.46

# ---------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.95999999999999
$status.Range("D4").Value = 13
$status.Range("E4").Value = -0.04
$status.Range("F4").Value = -0.04
$status.Range("G4").Value = 38.46

# ---------------------------------------------------------------------
# New trade row (#13) appended to "All Trades" and "MarketMaking" sheets
# ---------------------------------------------------------------------
$newRow = @(
    13,
    "2026-02-17",
    "07:53:44",
    "MarketMaking",
    "UP",
    0.73,
    0.7,
    "CLOSED",
    -4.1096,
    -0.03,
    99.95999999999999,
    0,
    0,
    0.6,
    "Normal spread capture: 19600 bps",
    "early_exit",
    0.15
)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Date/time columns (B, C) are stored as plain text in this workbook (not
    # real date serials), so force the cell format to Text before writing
    # them to avoid Excel auto-converting "2026-02-17" / "07:53:44" into
    # date/time serial numbers.
    $ws.Range("B14:C14").NumberFormat = "@"

    for ($i = 0; $i -lt $newRow.Length; $i++) {
        $ws.Cells.Item(14, $i + 1).Value = $newRow[$i]
    }

    # Restore the default "Normal" style so the new cells match the rest of
    # the sheet (which uses no explicit formatting).
    $ws.Range("B14:C14").Style = "Normal"
}
